$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209, pushing existing rows 209:337 down to 210:338.
$ws.Rows.Item(209).Insert()

# Populate the new row 209 with the added record.
$ws.Range("A209").Value = 4
$ws.Range("B209").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C209").Value = "Los Lagos"
$ws.Range("D209").Value = 44824
$ws.Range("E209").Value = 10
$ws.Range("F209").Value = 100112037
$ws.Range("G209").Value = "Cebollín"
$ws.Range("H209").Value = "Sin especificar"
$ws.Range("I209").Value = "Primera"
$ws.Range("J209").Value = 35
$ws.Range("K209").Value = 8500
$ws.Range("L209").Value = 8500
$ws.Range("M209").Value = 8500
$ws.Range("N209").Value = "$/paquete 36 unidades"
$ws.Range("O209").Value = "Región Metropolitana"
$ws.Range("P209").Value = 236
$ws.Range("Q209").Value = 36
$ws.Range("R209").Value = "Hortaliza"
